$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, [string]$text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") "55.141.18"
Set-TextValue $ws.Range("E2") "  -2.07%  "
Set-TextValue $ws.Range("D3") "2.355.37"
Set-TextValue $ws.Range("E3") "  -5.02%  "
Set-TextValue $ws.Range("E4") "  -0.03%  "
Set-TextValue $ws.Range("D5") "476.12"
Set-TextValue $ws.Range("E5") "  -2.43%  "
Set-TextValue $ws.Range("D6") "145.08"
Set-TextValue $ws.Range("E6") "  -1.15%  "
Set-TextValue $ws.Range("B7") "XRP"
Set-TextValue $ws.Range("C7") "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextValue $ws.Range("D7") "0.617"
Set-TextValue $ws.Range("E7") "  +20.81%  "
Set-TextValue $ws.Range("B8") "USDC"
Set-TextValue $ws.Range("C8") "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextValue $ws.Range("D8") "0.998"
Set-TextValue $ws.Range("E8") "  +0.14%  "
Set-TextValue $ws.Range("D9") "2.360.65"
Set-TextValue $ws.Range("E9") "  -5.07%  "
Set-TextValue $ws.Range("D10") "0.0973"
Set-TextValue $ws.Range("E10") "  +0.62%  "
Set-TextValue $ws.Range("E11") "  -5.85%  "
Set-TextValue $ws.Range("E12") "  -2.08%  "
Set-TextValue $ws.Range("E13") "  +1.08%  "
Set-TextValue $ws.Range("D14") "2.760.14"
Set-TextValue $ws.Range("E14") "  -5.27%  "
Set-TextValue $ws.Range("D15") "55.033.45"
Set-TextValue $ws.Range("E15") "  -2.25%  "
Set-TextValue $ws.Range("D16") "20.05"
Set-TextValue $ws.Range("E16") "  -4.77%  "
Set-TextValue $ws.Range("E17") "  -3.47%  "
Set-TextValue $ws.Range("D18") "2.355.78"
Set-TextValue $ws.Range("E18") "  -5.16%  "
Set-TextValue $ws.Range("E19") "  +1.75%  "
Set-TextValue $ws.Range("D20") "316.54"
Set-TextValue $ws.Range("E20") "  -0.28%  "
Set-TextValue $ws.Range("D21") "9.61"
Set-TextValue $ws.Range("E21") "  -4.45%  "
Set-TextValue $ws.Range("E22") "  +0.11%  "
Set-TextValue $ws.Range("D23") "5.63"
Set-TextValue $ws.Range("E23") "  -2.79%  "
Set-TextValue $ws.Range("D24") "56.94"
Set-TextValue $ws.Range("E24") "  -2.49%  "
Set-TextValue $ws.Range("D25") "1.00"
Set-TextValue $ws.Range("E25") "  +0.13%  "
Set-TextValue $ws.Range("D26") "0.394"
Set-TextValue $ws.Range("E26") "  -3.95%  "
Set-TextValue $ws.Range("E27") "  -4.76%  "
Set-TextValue $ws.Range("D28") "2.450.26"
Set-TextValue $ws.Range("E28") "  -5.26%  "
Set-TextValue $ws.Range("D29") "7.16"
Set-TextValue $ws.Range("E29") "  -5.95%  "
Set-TextValue $ws.Range("D30") "1.00"
Set-TextValue $ws.Range("E30") "  +0.15%  "
Set-TextValue $ws.Range("D31") "0.0₃0754"
Set-TextValue $ws.Range("E31") "  -4.85%  "
Set-TextValue $ws.Range("D32") "146.85"
Set-TextValue $ws.Range("E32") "  -1.90%  "
Set-TextValue $ws.Range("D33") "18.21"
Set-TextValue $ws.Range("E33") "  +0.37%  "
Set-TextValue $ws.Range("E34") "  -2.12%  "
Set-TextValue $ws.Range("D35") "5.08"
Set-TextValue $ws.Range("E35") "  -2.20%  "
Set-TextValue $ws.Range("E36") "  -3.96%  "
Set-TextValue $ws.Range("E37") "  -4.61%  "
Set-TextValue $ws.Range("D38") "0.814"
Set-TextValue $ws.Range("E38") "  -5.52%  "
Set-TextValue $ws.Range("D39") "33.73"
Set-TextValue $ws.Range("E39") "  -1.26%  "
Set-TextValue $ws.Range("B40") "Stellar"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D40") "0.0994"
Set-TextValue $ws.Range("E40") "  +7.42%  "
Set-TextValue $ws.Range("B41") "FirstDigitalUSD"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D41") "0.998"
Set-TextValue $ws.Range("E41") "  +0.44%  "
Set-TextValue $ws.Range("D42") "1.33"
Set-TextValue $ws.Range("E42") "  +0.38%  "
Set-TextValue $ws.Range("D43") "3.41"
Set-TextValue $ws.Range("E43") "  -2.67%  "
Set-TextValue $ws.Range("E44") "  -5.20%  "
Set-TextValue $ws.Range("D45") "0.0519"
Set-TextValue $ws.Range("E45") "  -6.63%  "
Set-TextValue $ws.Range("D46") "10.18"
Set-TextValue $ws.Range("E46") "  +0.06%  "
Set-TextValue $ws.Range("D47") "251.93"
Set-TextValue $ws.Range("E47") "  -2.87%  "
Set-TextValue $ws.Range("D48") "0.0221"
Set-TextValue $ws.Range("E48") "  -3.22%  "
Set-TextValue $ws.Range("D49") "4.36"
Set-TextValue $ws.Range("E49") "  -8.70%  "
Set-TextValue $ws.Range("D50") "16.73"
Set-TextValue $ws.Range("E50") "  -4.87%  "
Set-TextValue $ws.Range("D51") "1.780.10"
Set-TextValue $ws.Range("E51") "  -4.73%  "
